$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2566.1667
$ws.Range("I40").Value = 3180
$ws.Range("J40").Value = 2127.7144
$ws.Range("K40").Value = 3180
$ws.Range("L40").Value = 2127.7144
$ws.Range("M40").Value = -3005
$ws.Range("N40").Value = -2477.7144
$ws.Range("H92").Value = 43478604
$ws.Range("I92").Value = 58823868
$ws.Range("J92").Value = 358.5
$ws.Range("K92").Value = 58823868
$ws.Range("L92").Value = 358.5
$ws.Range("M92").Value = -58822620
$ws.Range("N92").Value = -2854.5
$ws.Range("H132").Value = 338333
$ws.Range("I132").Value = 666666
$ws.Range("K132").Value = 1999998
$ws.Range("M132").Value = -1997468

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1822.375
$ws.Range("I4").Value = 412.83334
$ws.Range("J4").Value = 6051
$ws.Range("K4").Value = 412.83334
$ws.Range("L4").Value = 6051
$ws.Range("M4").Value = -296.83334
$ws.Range("N4").Value = -6283
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").Value = $null
$ws.Range("H28").Value = 2409.2856
$ws.Range("I28").Value = 2409.2856
$ws.Range("K28").Value = 2409.2856
$ws.Range("M28").Value = -2217.2856
$ws.Range("H31").Value = 4943.6
$ws.Range("I31").Value = 4943.6
$ws.Range("K31").Value = 4943.6
$ws.Range("M31").Value = -4649.6
$ws.Range("H45").Value = 1718.2
$ws.Range("J45").Value = 1675
$ws.Range("L45").Value = 1675
$ws.Range("N45").Value = -2429
$ws.Range("H61").Value = 791.6667
$ws.Range("I61").Value = 791.6667
$ws.Range("K61").Value = 791.6667
$ws.Range("M61").Value = -579.6667
$ws.Range("H99").Value = 2409.2856
$ws.Range("I99").Value = 2409.2856
$ws.Range("K99").Value = 2409.2856
$ws.Range("M99").Value = 585.7143999999998
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -59178
$ws.Range("H128").Value = 50429
$ws.Range("J128").Value = 50429
$ws.Range("L128").Value = 50429
$ws.Range("N128").Value = -60389
$ws.Range("H136").Value = 791.6667
$ws.Range("I136").Value = 791.6667
$ws.Range("K136").Value = 2375.0001
$ws.Range("M136").Value = 174.9998999999998
$ws.Range("H140").Value = 59999.668
$ws.Range("J140").Value = 59999.668
$ws.Range("L140").Value = 59999.668
$ws.Range("N140").Value = -70359.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 602.8
$ws.Range("I22").Value = 549.36365
$ws.Range("K22").Value = 549.36365
$ws.Range("M22").Value = -376.36365
$ws.Range("H99").Value = 1642.6364
$ws.Range("I99").Value = 1556.9
$ws.Range("K99").Value = 1556.9
$ws.Range("M99").Value = -58.90000000000009
$ws.Range("H102").Value = 18101.285
$ws.Range("I102").Value = 18101.285
$ws.Range("K102").Value = 18101.285
$ws.Range("M102").Value = -14856.285
$ws.Range("H112").Value = 176666.33
$ws.Range("J112").Value = 176666.33
$ws.Range("L112").Value = 176666.33
$ws.Range("N112").Value = -179620.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4195.6875
$ws.Range("I31").Value = 1258.6
$ws.Range("K31").Value = 1258.6
$ws.Range("M31").Value = -963.5999999999999
$ws.Range("H34").Value = 4195.6875
$ws.Range("I34").Value = 1258.6
$ws.Range("K34").Value = 1258.6
$ws.Range("M34").Value = -1056.6
$ws.Range("H86").Value = 1000000000
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -999998877
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 1000000000
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").Value = $null
$ws.Range("H125").Value = 62500
$ws.Range("J125").Value = 62500
$ws.Range("L125").Value = 62500
$ws.Range("N125").Value = -67420

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 836.6667
$ws.Range("I60").Value = 836.6667
$ws.Range("K60").Value = 2510.0001
$ws.Range("M60").Value = -2259.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15172.546
$ws.Range("J49").Value = 24224.5
$ws.Range("L49").Value = 24224.5
$ws.Range("N49").Value = -24592.5
$ws.Range("I70").Value = 2497.5
$ws.Range("J70").Value = 250002750
$ws.Range("K70").Value = 2497.5
$ws.Range("L70").Value = 250002750
$ws.Range("M70").Value = -2227.5
$ws.Range("N70").Value = -250003290
$ws.Range("I73").Value = 2497.5
$ws.Range("J73").Value = 250002750
$ws.Range("K73").Value = 2497.5
$ws.Range("L73").Value = 250002750
$ws.Range("M73").Value = -1561.5
$ws.Range("N73").Value = -250004622
$ws.Range("H80").Value = 4235.8
$ws.Range("I80").Value = 4699.5
$ws.Range("J80").Value = 3926.6667
$ws.Range("K80").Value = 4699.5
$ws.Range("L80").Value = 3926.6667
$ws.Range("M80").Value = -3701.5
$ws.Range("N80").Value = -5922.6667
$ws.Range("H83").Value = 4235.8
$ws.Range("I83").Value = 4699.5
$ws.Range("J83").Value = 3926.6667
$ws.Range("K83").Value = 23497.5
$ws.Range("L83").Value = 19633.3335
$ws.Range("M83").Value = -18505.5
$ws.Range("N83").Value = -29617.3335
$ws.Range("H134").Value = 46663
$ws.Range("J134").Value = 46663
$ws.Range("L134").Value = 139989
$ws.Range("N134").Value = -145059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 227.14285
$ws.Range("H27").Value = 227.14285
$ws.Range("H55").Value = 1503.7858
$ws.Range("J55").Value = 1913.1666
$ws.Range("L55").Value = 1913.1666
$ws.Range("N55").Value = -2259.1666
$ws.Range("H93").Value = 47621450
$ws.Range("I93").Value = 83335810
$ws.Range("K93").Value = 83335810
$ws.Range("M93").Value = -83334562
$ws.Range("H120").Value = 7349
$ws.Range("J120").Value = 7349
$ws.Range("L120").Value = 7349
$ws.Range("N120").Value = -17025
$ws.Range("H132").Value = 2255.8
$ws.Range("I132").Value = 1944.75
$ws.Range("K132").Value = 5834.25
$ws.Range("M132").Value = -3304.25
$ws.Range("H135").Value = 129999.5
$ws.Range("J135").Value = 129999.5
$ws.Range("L135").Value = 129999.5
$ws.Range("N135").Value = -140139.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 3500
$ws.Range("I29").Value = 2000
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 2000
$ws.Range("L29").Value = 5000
$ws.Range("M29").Value = -1710
$ws.Range("N29").Value = -5580
$ws.Range("H100").Value = 2360.1875
$ws.Range("I100").Value = 2334.625
$ws.Range("J100").Value = 2385.75
$ws.Range("K100").Value = 4669.25
$ws.Range("L100").Value = 4771.5
$ws.Range("M100").Value = -4128.25
$ws.Range("N100").Value = -5853.5
$ws.Range("H113").Value = 601.4286
$ws.Range("I113").Value = 451.14285
$ws.Range("K113").Value = 1353.42855
$ws.Range("M113").Value = 816.5714499999999
$ws.Range("H136").Value = 1805.5555
$ws.Range("I136").Value = 1805.5555
$ws.Range("K136").Value = 5416.666499999999
$ws.Range("M136").Value = -2866.666499999999
